# Apply odds updates for week games (2025-03-19) per commit "Atualizando o arquivo XLSX"
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("H3").Value = 3.2
$ws.Range("J3").Value = 2.5
$ws.Range("K3").Value = 2
$ws.Range("M3").Value = 1.08
$ws.Range("N3").Value = 8
$ws.Range("O3").Value = 1.44
$ws.Range("P3").Value = 2.63
$ws.Range("Q3").Value = 1.8
$ws.Range("R3").Value = 2.05
$ws.Range("S3").Value = 2.35
$ws.Range("T3").Value = 1.57
$ws.Range("U3").Value = 3.9
$ws.Range("V3").Value = 1.25
$ws.Range("W3").Value = 4.5
$ws.Range("X3").Value = 1.18
$ws.Range("Y3").Value = 1.53
$ws.Range("Z3").Value = 2.38
$ws.Range("AA3").Value = 2.2
$ws.Range("AB3").Value = 1.62
$ws.Range("AD3").Value = 7.5
$ws.Range("AE3").Value = 9
$ws.Range("AG3").Value = 17
$ws.Range("AH3").Value = 34
$ws.Range("AI3").Value = 7
$ws.Range("AK3").Value = 19
$ws.Range("AL3").Value = 67
$ws.Range("AP3").Value = 17
$ws.Range("G4").Value = 1.7
$ws.Range("H4").Value = 3.7
$ws.Range("I4").Value = 5
$ws.Range("J4").Value = 2.38
$ws.Range("K4").Value = 2.05
$ws.Range("L4").Value = 5.5
$ws.Range("S4").Value = 2.2
$ws.Range("T4").Value = 1.65
$ws.Range("U4").Value = 3.35
$ws.Range("Y4").Value = 1.47
$ws.Range("Z4").Value = 2.5
$ws.Range("AF4").Value = 13
$ws.Range("AG4").Value = 17
$ws.Range("AI4").Value = 8
$ws.Range("AJ4").Value = 7
$ws.Range("AL4").Value = 67
$ws.Range("AO4").Value = 23
$ws.Range("G5").Value = 4.2
$ws.Range("H5").Value = 3.15
$ws.Range("I5").Value = 1.87
$ws.Range("J5").Value = 4.5
$ws.Range("K5").Value = 1.98
$ws.Range("L5").Value = 2.55
$ws.Range("P5").Value = 3
$ws.Range("S5").Value = 1.85
$ws.Range("T5").Value = 1.75
$ws.Range("W5").Value = 2.92
$ws.Range("X5").Value = 1.3
$ws.Range("Z5").Value = 2.45
$ws.Range("AA5").Value = 1.7
$ws.Range("AB5").Value = 1.93
$ws.Range("AC5").Value = 12.5
$ws.Range("AD5").Value = 26
$ws.Range("AE5").Value = 13.5
$ws.Range("AF5").Value = 75
$ws.Range("AG5").Value = 37
$ws.Range("AH5").Value = 37
$ws.Range("AI5").Value = 9
$ws.Range("AJ5").Value = 6.2
$ws.Range("AK5").Value = 13
$ws.Range("AL5").Value = 55
$ws.Range("AM5").Value = 400
$ws.Range("AN5").Value = 6.9
$ws.Range("AO5").Value = 8.75
$ws.Range("AQ5").Value = 16.5
$ws.Range("AR5").Value = 15
